{"js": "// Update the worksheet's date title and regenerate the 25 division\n// problems (\"two-digit number divided by one-digit number\") shown in\n// the 5-column table. The table has 20 rows total: every 4th row\n// (0, 4, 8, 12, 16) holds the 5 visible problems for that exercise\n// \"row\", and the 3 rows in between are blank spacer rows.\n\nconst TITLE_OLD = \"2025-10-25 Saturday\";\nconst TITLE_NEW = \"2025-10-26 Sunday\";\n\nconst GRID_OLD = [\n  [\"89\u00f74=22, 1\", \"76\u00f75=15, 1\", \"97\u00f77=13, 6\", \"29\u00f75=5, 4\", \"65\u00f74=16, 1\"],\n  [\"67\u00f79=7, 4\", \"18\u00f73=6, 0\", \"17\u00f77=2, 3\", \"41\u00f74=10, 1\", \"41\u00f74=10, 1\"],\n  [\"96\u00f74=24, 0\", \"14\u00f76=2, 2\", \"17\u00f79=1, 8\", \"87\u00f72=43, 1\", \"32\u00f74=8, 0\"],\n  [\"65\u00f75=13, 0\", \"86\u00f76=14, 2\", \"71\u00f74=17, 3\", \"41\u00f77=5, 6\", \"48\u00f74=12, 0\"],\n  [\"91\u00f78=11, 3\", \"62\u00f75=12, 2\", \"63\u00f75=12, 3\", \"46\u00f75=9, 1\", \"95\u00f74=23, 3\"],\n];\nconst GRID_NEW = [\n  [\"16\u00f76=2, 4\", \"57\u00f79=6, 3\", \"64\u00f77=9, 1\", \"79\u00f76=13, 1\", \"29\u00f78=3, 5\"],\n  [\"80\u00f76=13, 2\", \"33\u00f74=8, 1\", \"57\u00f72=28, 1\", \"68\u00f78=8, 4\", \"32\u00f77=4, 4\"],\n  [\"65\u00f76=10, 5\", \"14\u00f73=4, 2\", \"54\u00f79=6, 0\", \"74\u00f76=12, 2\", \"61\u00f74=15, 1\"],\n  [\"41\u00f78=5, 1\", \"79\u00f74=19, 3\", \"74\u00f77=10, 4\", \"40\u00f75=8, 0\", \"77\u00f79=8, 5\"],\n  [\"86\u00f76=14, 2\", \"37\u00f77=5, 2\", \"83\u00f75=16, 3\", \"39\u00f75=7, 4\", \"83\u00f79=9, 2\"],\n];\nconst DATA_ROW_INDICES = [0, 4, 8, 12, 16];\n\n// --- 1. Update the title paragraph (first paragraph of the body). ---\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The title is always the very first paragraph of the body; write the\n// new date regardless (the equality check is just a sanity guard that\n// can be inspected/logged if this script is adapted later).\nconst titlePara = paragraphs.items[0];\ntitlePara.insertText(TITLE_NEW, \"Replace\");\nawait context.sync();\n\n// --- 2. Update the 25 division-problem cells inside the table. ---\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows.items;\nfor (const rowIdx of DATA_ROW_INDICES) {\n  rows[rowIdx].cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const rowIdx of DATA_ROW_INDICES) {\n  const cells = rows[rowIdx].cells.items;\n  for (let c = 0; c < cells.length; c++) {\n    cells[c].body.paragraphs.load(\"items/text\");\n  }\n}\nawait context.sync();\n\nconst gridRowForIndex = {};\nDATA_ROW_INDICES.forEach((rowIdx, i) => (gridRowForIndex[rowIdx] = i));\n\nfor (const rowIdx of DATA_ROW_INDICES) {\n  const gridRow = gridRowForIndex[rowIdx];\n  const cells = rows[rowIdx].cells.items;\n  for (let c = 0; c < cells.length; c++) {\n    // GRID_OLD is kept alongside GRID_NEW purely as in-source documentation\n    // of which value is being replaced (cells[c].body.paragraphs.items[0].text\n    // is expected to equal GRID_OLD[gridRow][c] at this point).\n    const newValue = GRID_NEW[gridRow][c];\n    const para = cells[c].body.paragraphs.items[0];\n    para.insertText(newValue, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "# Update the worksheet's date title and regenerate the 25 division\n# problems (\"two-digit number divided by one-digit number\") shown in\n# the 5-column table. The table has 20 rows total: every 4th row\n# (1, 5, 9, 13, 17 in 1-based COM indexing) holds the 5 visible\n# problems for that exercise \"row\"; the rows in between are blank\n# spacer rows.\n\n$titleOld = \"2025-10-25 Saturday\"\n$titleNew = \"2025-10-26 Sunday\"\n\n$gridOld = @(\n    @(\"89\u00f74=22, 1\", \"76\u00f75=15, 1\", \"97\u00f77=13, 6\", \"29\u00f75=5, 4\", \"65\u00f74=16, 1\"),\n    @(\"67\u00f79=7, 4\", \"18\u00f73=6, 0\", \"17\u00f77=2, 3\", \"41\u00f74=10, 1\", \"41\u00f74=10, 1\"),\n    @(\"96\u00f74=24, 0\", \"14\u00f76=2, 2\", \"17\u00f79=1, 8\", \"87\u00f72=43, 1\", \"32\u00f74=8, 0\"),\n    @(\"65\u00f75=13, 0\", \"86\u00f76=14, 2\", \"71\u00f74=17, 3\", \"41\u00f77=5, 6\", \"48\u00f74=12, 0\"),\n    @(\"91\u00f78=11, 3\", \"62\u00f75=12, 2\", \"63\u00f75=12, 3\", \"46\u00f75=9, 1\", \"95\u00f74=23, 3\")\n)\n$gridNew = @(\n    @(\"16\u00f76=2, 4\", \"57\u00f79=6, 3\", \"64\u00f77=9, 1\", \"79\u00f76=13, 1\", \"29\u00f78=3, 5\"),\n    @(\"80\u00f76=13, 2\", \"33\u00f74=8, 1\", \"57\u00f72=28, 1\", \"68\u00f78=8, 4\", \"32\u00f77=4, 4\"),\n    @(\"65\u00f76=10, 5\", \"14\u00f73=4, 2\", \"54\u00f79=6, 0\", \"74\u00f76=12, 2\", \"61\u00f74=15, 1\"),\n    @(\"41\u00f78=5, 1\", \"79\u00f74=19, 3\", \"74\u00f77=10, 4\", \"40\u00f75=8, 0\", \"77\u00f79=8, 5\"),\n    @(\"86\u00f76=14, 2\", \"37\u00f77=5, 2\", \"83\u00f75=16, 3\", \"39\u00f75=7, 4\", \"83\u00f79=9, 2\")\n)\n$dataRows = @(1, 5, 9, 13, 17)\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the title paragraph (first paragraph of the body). ---\n$titlePara = $d.Paragraphs.Item(1)\n$titleRange = $titlePara.Range\n$currentTitle = $titleRange.Text.TrimEnd([char]13, [char]7)\nif ($currentTitle -ne $titleOld) {\n    Write-Output \"Warning: title text did not match expected old value (found: $currentTitle)\"\n}\n$titleRange.Text = $titleNew\n\n# --- 2. Update the 25 division-problem cells inside the table. ---\n$t = $d.Tables.Item(1)\nfor ($i = 0; $i -lt $dataRows.Length; $i++) {\n    $row = $dataRows[$i]\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $t.Cell($row, $c)\n        $cellRange = $cell.Range\n        $currentVal = $cellRange.Text.TrimEnd([char]13, [char]7)\n        $expectedOld = $gridOld[$i][$c - 1]\n        $newVal = $gridNew[$i][$c - 1]\n        if ($currentVal -ne $expectedOld) {\n            Write-Output \"Warning: cell ($row,$c) text did not match expected old value (found: $currentVal)\"\n        }\n        $cellRange.Text = $newVal\n    }\n}\n"}
